# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Rule row 11 (B11) on the "Rules" sheet changes from the text "R40" to
# the text "1". Format the cell as Text first so Excel stores the new
# value as a (shared) string rather than re-interpreting "1" as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
